$d = $word.ActiveDocument

# 1. Update the title line text:
#    "Products and Services Exchange Network:"
#    -> "Products And Services Community Exchange Network:"
$oldTitle = "Products and Services Exchange Network:"
$newTitle = "Products And Services Community Exchange Network:"
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $oldTitle) {
        $p.Range.Text = $newTitle
        break
    }
}

# 2. After the "Domains: data, schema and behavior ..." paragraph, insert a
#    blank paragraph followed by a new paragraph containing the extra
#    sentence about general purpose business domains helper tools.
$targetText = "Domains: data, schema and behavior of business applications (ERP, CRM, BI, SCM, HMS, etc.)."
$newSentence = "General purpose business domains problem resolution / tasks, goals accomplishment helper tools."
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $targetText) {
        $p.Range.InsertParagraphAfter()
        $blankPara = $d.Paragraphs.Item($i + 1)
        $blankPara.Range.InsertParagraphAfter()
        $textPara = $d.Paragraphs.Item($i + 2)
        $textPara.Range.Text = $newSentence
        break
    }
}
